# DevLog.xlsx update
# - Insert a new milestone row at row 9 (AVL tree work: "Not a milestone -
#   Implemented all functions aside from test_avl in avl_tree.c...")
# - Remove the blank spacer row that used to sit right below the last real
#   entry (old row 14)
# - Append a new blank spacer row at the very bottom of the table
# - Update the sheet's scroll position / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a brand-new row at row 9 and give it the same look & feel as
#    the row directly below it (which used to be row 9, the long wrapped
#    entry), then overwrite the text/values for the AVL milestone.
# ---------------------------------------------------------------------
$ws.Rows("9:9").Insert()

$ws.Range("A10:G10").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# (values are written in the same order the matching shared-string
# entries appear in the workbook, i.e. Milestone -> Test Functions ->
# Resources -> Functions -> Implementation notes)
$ws.Range("B9").Value = "Not a milestone - Implemented all functions aside from test_avl in avl_tree.c. Also created a header file for it."
$ws.Range("C9").Value = 45771
$ws.Range("F9").Value = "test_avl - not fully implemented. Only holdds adhoc tests"
$ws.Range("E9").Value = "The tutorials, specifically the week 5 tutorial on BSTs, were used heavily. Code was copied from it, and functions were derivative of it. ChatGPT was used for gaining a slightly better understanding of how balancing the tree could be implemented but for the most part the balancing and insertion methods were derived from my notes on the weekly content."
$ws.Range("D9").Value = "create_avl`nfind_avl_node`nfind_avl`nfind_left_right_height`nget_avl_node_height`nleft_rotate`nright_rotate`nbalance_tree`ninsert_avl_node`ninsert_avl`nprint_avl_node`nprint_avl`ntest_avl"
$ws.Range("G9").Value = "test_avl has been implemented with some adhoc tests that served purely to see the basic functionality of the functions. These are by no means exhaustive and do not at all match the format seen in other testing functions. This file also has more code from the tutorials than any other file. Both the linked list and record files were largely implemented with very little tutorial code and were largely custom. This file also uses a similar format of having wrapper functions and a wrapper struct as the BST file did in the week 5 tutorial."

# Let the new row grow tall enough to show the long note in G9.
$ws.Rows("9:9").RowHeight = 195

# ---------------------------------------------------------------------
# 2. Old row 14 was an empty spacer row directly under the table (now at
#    row 15 since we inserted a row above). Delete it - it is no longer
#    needed because a new spacer row is added at the bottom instead.
# ---------------------------------------------------------------------
$ws.Rows("15:15").Delete()

# ---------------------------------------------------------------------
# 3. Add a fresh blank spacer row at the bottom of the sheet (row 53),
#    matching the formatting the old last row (52) used to have, and
#    restyle row 52 to match the rest of the blank filler rows above it.
# ---------------------------------------------------------------------
$ws.Range("B52:G52").Copy()
$ws.Range("B53:G53").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B51:G51").Copy()
$ws.Range("B52:G52").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. Restore the on-screen selection/scroll position.
# ---------------------------------------------------------------------
$ws.Range("F9").Select()
